# Add a new worksheet "EQ_CCP" between the existing "TG_CCP" and "TG_OOS"
# sheets, and populate it with the new test-data grid (Variable Name /
# Valid / Invalid columns), matching the authored diff.

$wb = $excel.ActiveWorkbook

$tgCcp = $wb.Worksheets.Item("TG_CCP")

# Insert the new sheet right after TG_CCP (i.e. before TG_OOS).
$eqCcp = $wb.Worksheets.Add($null, $tgCcp)
$eqCcp.Name = "EQ_CCP"

# Header row (row 5), left to right.
$eqCcp.Cells.Item(5, 3).Value = "Sr. No."
$eqCcp.Cells.Item(5, 4).Value = "Variable Name"
$eqCcp.Cells.Item(5, 5).Value = "Valid"
$eqCcp.Cells.Item(5, 6).Value = "Invalid"

# "Valid" column for the two data rows.
$eqCcp.Cells.Item(6, 5).Value = "(a-z),(0-9)"
$eqCcp.Cells.Item(7, 5).Value = '(a-z),(0-9),"@",".","_"'

# "Variable Name" column for the two data rows.
$eqCcp.Cells.Item(7, 4).Value = "Email"
$eqCcp.Cells.Item(6, 4).Value = "User Name"

# "Invalid" column for the two data rows.
$eqCcp.Cells.Item(6, 6).Value = '{{()}]!@#$%^&|,>.*+/-?\~<`:";'''
$eqCcp.Cells.Item(7, 6).Value = '{{()}]!#$%^&|,>.*+/-?\~<`:";'''

# "Sr. No." column, rows 6-15 (1..10).
$eqCcp.Cells.Item(6, 3).Value = 1
$eqCcp.Cells.Item(7, 3).Value = 2
$eqCcp.Cells.Item(8, 3).Value = 3
$eqCcp.Cells.Item(9, 3).Value = 4
$eqCcp.Cells.Item(10, 3).Value = 5
$eqCcp.Cells.Item(11, 3).Value = 6
$eqCcp.Cells.Item(12, 3).Value = 7
$eqCcp.Cells.Item(13, 3).Value = 8
$eqCcp.Cells.Item(14, 3).Value = 9
$eqCcp.Cells.Item(15, 3).Value = 10

# Column widths roughly matching the authored sheet.
$eqCcp.Columns.Item(4).ColumnWidth = 15.67
$eqCcp.Columns.Item(5).ColumnWidth = 24.67
$eqCcp.Columns.Item(6).ColumnWidth = 36
$eqCcp.Columns.Item(7).ColumnWidth = 13.5

# Match the authored file's final selection on the new sheet.
$eqCcp.Range("F7").Select()
